$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells: "_old" -> "_FV2404", "_new" -> "_FV2410"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if (($val -ne $null) -and ($val -is [string])) {
        if ($val.EndsWith("_old")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2404"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2410"
        }
    }
}

# 2. Turn the range into an Excel Table (ListObject)
$listObject = $ws.ListObjects.Add(1, $ws.Range("A1:U85"), 0, 1)
$listObject.Name = "Table1"

# 3. Freeze the header row (pane split after row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
